# "Add area and efficiency"
#  - Recalculate/update the PV area figure for "Tower A" (row 20, column B)
#  - Leave the cursor/selection parked on B20 (the cell that was edited)
#  - Set up the sheet's page setup (paper size / orientation) for printing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B20 ("Tower A" area) to the corrected value
$ws.Range("B20").Value = 166.55758399999999

# Scroll the view near row 11 and move the active selection to B20,
# matching where the workbook was left after the edit
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("B20").Select()

# Configure the page for printing: A4, portrait
$ps = $ws.PageSetup
$ps.PaperSize = [Microsoft.Office.Interop.Excel.XlPaperSize]::xlPaperA4
$ps.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait
